$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 — update B2/C2 (unchanged text, new shared-string slot), D2, and the
# numeric weight/specificity columns with the recomputed values.
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Ccl4"
$ws.Range("C2").Value = "Ccr5"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 335.6132
$ws.Range("H2").Value = 1006.8396
$ws.Range("I2").Value = 0.80464917790985
$ws.Range("J2").Value = 0.80464917790985
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 4.021407666666667
$ws.Range("N2").Value = 12.064223
$ws.Range("O2").Value = 0.06269882270324605
$ws.Range("P2").Value = 0.06269882270324605
$ws.Range("Q2").Value = 1349.637495514533
$ws.Range("R2").Value = 12146.7374596308
$ws.Range("S2").Value = 0.05045055614408237
$ws.Range("T2").Value = 0.05045055614408237

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Ccl4"
$ws.Range("C3").Value = "Ccr5"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 335.6132
$ws.Range("H3").Value = 1006.8396
$ws.Range("I3").Value = 0.80464917790985
$ws.Range("J3").Value = 0.80464917790985
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.4010506666666667
$ws.Range("N3").Value = 1.203152
$ws.Range("O3").Value = 0.006252886235031953
$ws.Range("P3").Value = 0.006252886235031953
$ws.Range("Q3").Value = 134.5978976021333
$ws.Range("R3").Value = 1211.3810784192
$ws.Range("S3").Value = 0.005031379768582278
$ws.Range("T3").Value = 0.005031379768582278

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Ccl4"
$ws.Range("C4").Value = "Ccr5"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 335.6132
$ws.Range("H4").Value = 1006.8396
$ws.Range("I4").Value = 0.80464917790985
$ws.Range("J4").Value = 0.80464917790985
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 59.71602933333333
$ws.Range("N4").Value = 179.148088
$ws.Range("O4").Value = 0.931048291061722
$ws.Range("P4").Value = 0.931048291061722
$ws.Range("Q4").Value = 20041.48769585387
$ws.Range("R4").Value = 180373.3892626848
$ws.Range("S4").Value = 0.7491672419971854
$ws.Range("T4").Value = 0.7491672419971854

# Row 5
$ws.Range("A5").Value = "M2"
$ws.Range("B5").Value = "Ccl4"
$ws.Range("C5").Value = "Ccr5"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 81.47937800000001
$ws.Range("H5").Value = 244.438134
$ws.Range("I5").Value = 0.19535082209015
$ws.Range("J5").Value = 0.19535082209015
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 4.021407666666667
$ws.Range("N5").Value = 12.064223
$ws.Range("O5").Value = 0.06269882270324605
$ws.Range("P5").Value = 0.06269882270324605
$ws.Range("Q5").Value = 327.6617953644314
$ws.Range("R5").Value = 2948.956158279882
$ws.Range("S5").Value = 0.01224826655916368
$ws.Range("T5").Value = 0.01224826655916367

# Row 6
$ws.Range("A6").Value = "M2"
$ws.Range("B6").Value = "Ccl4"
$ws.Range("C6").Value = "Ccr5"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 81.47937800000001
$ws.Range("H6").Value = 244.438134
$ws.Range("I6").Value = 0.19535082209015
$ws.Range("J6").Value = 0.19535082209015
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.4010506666666667
$ws.Range("N6").Value = 1.203152
$ws.Range("O6").Value = 0.006252886235031953
$ws.Range("P6").Value = 0.006252886235031953
$ws.Range("Q6").Value = 32.67735886648533
$ws.Range("R6").Value = 294.096229798368
$ws.Range("S6").Value = 0.001221506466449675
$ws.Range("T6").Value = 0.001221506466449675

# Row 7
$ws.Range("A7").Value = "M2"
$ws.Range("B7").Value = "Ccl4"
$ws.Range("C7").Value = "Ccr5"
$ws.Range("D7").Value = "M2"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 81.47937800000001
$ws.Range("H7").Value = 244.438134
$ws.Range("I7").Value = 0.19535082209015
$ws.Range("J7").Value = 0.19535082209015
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 59.71602933333333
$ws.Range("N7").Value = 179.148088
$ws.Range("O7").Value = 0.931048291061722
$ws.Range("P7").Value = 0.931048291061722
$ws.Range("Q7").Value = 4865.624926709756
$ws.Range("R7").Value = 43790.6243403878
$ws.Range("S7").Value = 0.1818810490645367
$ws.Range("T7").Value = 0.1818810490645366
